$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.063.24"
$ws.Range("E2").Value = "  +4.70%  "
$ws.Range("D3").Value = "2.233.02"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'245.87"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "'75.59"
$ws.Range("E7").Value = "  +8.06%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.613"
$ws.Range("E9").Value = "  +6.30%  "
$ws.Range("D10").Value = "'41.16"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "'55.46"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'6.98"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "2.569.48"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").Value = "'14.72"
$ws.Range("E16").Value = "  +5.85%  "
$ws.Range("D17").Value = "2.241.46"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("D18").Value = "'0.810"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "42.941.79"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "'10.46"
$ws.Range("E23").Value = "  +5.06%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.21"
$ws.Range("E24").Value = "  +13.03%  "
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").Value = "'228.93"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'174.34"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.18"
$ws.Range("E31").Value = "  +20.92%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").Value = "'20.37"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("E35").Value = "  +4.22%  "
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").Value = "'0.112"
$ws.Range("E37").Value = "  +7.54%  "
$ws.Range("D38").Value = "'4.37"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "'0.0334"
$ws.Range("E39").Value = "  +17.14%  "
$ws.Range("D40").Value = "'13.06"
$ws.Range("E40").Value = "  +5.32%  "
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("D42").Value = "'5.59"
$ws.Range("E42").Value = "  +3.34%  "
$ws.Range("E43").Value = "  +5.32%  "
$ws.Range("D44").Value = "'60.13"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "'105.41"
$ws.Range("E45").Value = "  +7.33%  "
$ws.Range("D46").Value = "'8.57"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.33"
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").Value = "'0.442"
$ws.Range("E49").Value = "  +19.86%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'1.11"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("E51").Value = "  +2.07%  "
